# Update gh-pages to output generated at 456a3b4
# Increment "想去人数" (F column) counts for a few events across sheets.

$wb = $excel.ActiveWorkbook

$wsExpo = $wb.Worksheets.Item("展览")
$wsExpo.Range("F3").Value = 432

$wsShow = $wb.Worksheets.Item("演出")
$wsShow.Range("F2").Value = 70

$wsAll = $wb.Worksheets.Item("全部类型")
$wsAll.Range("F3").Value = 70
$wsAll.Range("F4").Value = 432
